$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 11
$ws.Range("H11").Value = 2781054.8
$ws.Range("I11").Value = 2781054.8
$ws.Range("K11").Value = 2781054.8
$ws.Range("M11").Value = -2780914.8
# row 46
$ws.Range("H46").Value = 18527456
$ws.Range("I46").Value = 37050796
$ws.Range("J46").Value = 4115.3335
$ws.Range("K46").Value = 111152388
$ws.Range("L46").Value = 12346.0005
$ws.Range("M46").Value = -111152269
$ws.Range("N46").Value = -12584.0005
# row 60
$ws.Range("H60").Value = 18527456
$ws.Range("I60").Value = 37050796
$ws.Range("J60").Value = 4115.3335
$ws.Range("K60").Value = 111152388
$ws.Range("L60").Value = 12346.0005
$ws.Range("M60").Value = -111151904
$ws.Range("N60").Value = -13314.0005
# row 64
$ws.Range("H64").Value = 3300
$ws.Range("I64").Value = 3293.3333
$ws.Range("K64").Value = 3293.3333
$ws.Range("M64").Value = -3045.3333
# row 67
$ws.Range("H67").Value = 3300
$ws.Range("I67").Value = 3293.3333
$ws.Range("K67").Value = 3293.3333
$ws.Range("M67").Value = -2435.3333
# row 80
$ws.Range("H80").Value = 15874285
$ws.Range("I80").Value = 55556492
$ws.Range("J80").Value = 1402.5333
$ws.Range("K80").Value = 166669476
$ws.Range("L80").Value = 4207.5999
$ws.Range("M80").Value = -166668478
$ws.Range("N80").Value = -6203.5999
# row 83
$ws.Range("H83").Value = 15874285
$ws.Range("I83").Value = 55556492
$ws.Range("J83").Value = 1402.5333
$ws.Range("K83").Value = 500008428
$ws.Range("L83").Value = 12622.7997
$ws.Range("M83").Value = -500003436
$ws.Range("N83").Value = -22606.7997
# row 132
$ws.Range("H132").Value = 26403.527
$ws.Range("I132").Value = 3713.8518
$ws.Range("J132").Value = 94472.55499999999
$ws.Range("K132").Value = 11141.5554
$ws.Range("L132").Value = 283417.665
$ws.Range("M132").Value = -8611.555399999999
$ws.Range("N132").Value = -288477.665

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 11091.947
$ws.Range("I32").Value = 9912.294
$ws.Range("J32").Value = 21119
$ws.Range("K32").Value = 9912.294
$ws.Range("L32").Value = 21119
$ws.Range("M32").Value = -9625.294
$ws.Range("N32").Value = -21693
# row 61
$ws.Range("H61").Value = 2491.5833
$ws.Range("I61").Value = 1780.6316
$ws.Range("J61").Value = 3286.1765
$ws.Range("K61").Value = 1780.6316
$ws.Range("L61").Value = 3286.1765
$ws.Range("M61").Value = -1568.6316
$ws.Range("N61").Value = -3710.1765
# row 97
$ws.Range("H97").Value = 3000
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 3000
$ws.Range("N97").Value = -3992
$ws.Range("M97").ClearContents()
# row 102
$ws.Range("H102").Value = 32656.5
$ws.Range("I102").Value = 4422
$ws.Range("K102").Value = 4422
$ws.Range("M102").Value = -2800
# row 122
$ws.Range("H122").Value = 1927.3334
$ws.Range("I122").Value = 1830.2307
$ws.Range("J122").Value = 2179.8
$ws.Range("K122").Value = 5490.6921
$ws.Range("L122").Value = 6539.400000000001
$ws.Range("M122").Value = -3040.6921
$ws.Range("N122").Value = -11439.4
# row 132
$ws.Range("H132").Value = 10206427
$ws.Range("I132").Value = 16130505
$ws.Range("J132").Value = 3847.7778
$ws.Range("K132").Value = 48391515
$ws.Range("L132").Value = 11543.3334
$ws.Range("M132").Value = -48388985
$ws.Range("N132").Value = -16603.3334
# row 136
$ws.Range("H136").Value = 2491.5833
$ws.Range("I136").Value = 1780.6316
$ws.Range("J136").Value = 3286.1765
$ws.Range("K136").Value = 5341.8948
$ws.Range("L136").Value = 9858.529500000001
$ws.Range("M136").Value = -2791.8948
$ws.Range("N136").Value = -14958.5295

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 20
$ws.Range("H20").Value = 5576.393
$ws.Range("I20").Value = 1571.3572
$ws.Range("J20").Value = 9581.429
$ws.Range("K20").Value = 1571.3572
$ws.Range("L20").Value = 9581.429
$ws.Range("M20").Value = -1324.3572
$ws.Range("N20").Value = -10075.429
# row 75
$ws.Range("H75").Value = 15020.083
$ws.Range("I75").Value = 5060.25
$ws.Range("J75").Value = 20000
$ws.Range("K75").Value = 5060.25
$ws.Range("L75").Value = 20000
$ws.Range("N75").Value = -21872
$ws.Range("M75").Value = -4124.25
# row 78
$ws.Range("H78").Value = 15020.083
$ws.Range("I78").Value = 5060.25
$ws.Range("J78").Value = 20000
$ws.Range("K78").Value = 15180.75
$ws.Range("L78").Value = 60000
$ws.Range("N78").Value = -69360
$ws.Range("M78").Value = -10500.75
# row 86
$ws.Range("H86").Value = 3500
$ws.Range("I86").Value = 3500
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3500
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -2377
$ws.Range("N86").ClearContents()
# row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
# row 89
$ws.Range("H89").Value = 3500
$ws.Range("I89").Value = 3500
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 17500
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -11884
$ws.Range("N89").Value = -11884
# row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
# row 94
$ws.Range("H94").Value = 2062.0454
$ws.Range("I94").Value = 2013.2142
$ws.Range("J94").Value = 2147.5
$ws.Range("K94").Value = 2013.2142
$ws.Range("L94").Value = 2147.5
$ws.Range("M94").Value = -1562.2142
$ws.Range("N94").Value = -3049.5
# row 99
$ws.Range("H99").Value = 2095.641
$ws.Range("I99").Value = 2021.4706
$ws.Range("J99").Value = 2600
$ws.Range("K99").Value = 2021.4706
$ws.Range("L99").Value = 2600
$ws.Range("M99").Value = -523.4706000000001
$ws.Range("N99").Value = -5596
# row 105
$ws.Range("H105").Value = 3534.2
$ws.Range("I105").Value = 2206
$ws.Range("J105").Value = 4249.385
$ws.Range("K105").Value = 2206
$ws.Range("L105").Value = 4249.385
$ws.Range("M105").Value = -459
$ws.Range("N105").Value = -7743.385

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 6949993
$ws.Range("I31").Value = 1995.3572
$ws.Range("K31").Value = 1995.3572
$ws.Range("M31").Value = -1700.3572
# row 34
$ws.Range("H34").Value = 6949993
$ws.Range("I34").Value = 1995.3572
$ws.Range("K34").Value = 1995.3572
$ws.Range("M34").Value = -1793.3572
# row 122
$ws.Range("H122").Value = 184643
$ws.Range("I122").Value = 315387.44
$ws.Range("J122").Value = 1600.8
$ws.Range("K122").Value = 946162.3200000001
$ws.Range("L122").Value = 4802.4
$ws.Range("M122").Value = -943712.3200000001
$ws.Range("N122").Value = -9702.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 38
$ws.Range("H38").Value = 125358.625
$ws.Range("J38").Value = 200470.2
$ws.Range("L38").Value = 601410.6000000001
$ws.Range("N38").Value = -602104.6000000001
# row 113
$ws.Range("H113").Value = 2760.149
$ws.Range("I113").Value = 4096.393
$ws.Range("J113").Value = 790.9474
$ws.Range("K113").Value = 12289.179
$ws.Range("L113").Value = 2372.8422
$ws.Range("M113").Value = -10119.179
$ws.Range("N113").Value = -6712.8422

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 5348
$ws.Range("I70").Value = 5371.4287
$ws.Range("J70").Value = 5225
$ws.Range("K70").Value = 5371.4287
$ws.Range("L70").Value = 5225
$ws.Range("M70").Value = -5101.4287
$ws.Range("N70").Value = -5765
# row 73
$ws.Range("H73").Value = 5348
$ws.Range("I73").Value = 5371.4287
$ws.Range("J73").Value = 5225
$ws.Range("K73").Value = 5371.4287
$ws.Range("L73").Value = 5225
$ws.Range("M73").Value = -4435.4287
$ws.Range("N73").Value = -7097
# row 97
$ws.Range("H97").Value = 2504.081
$ws.Range("I97").Value = 1897.0968
$ws.Range("J97").Value = 5640.1665
$ws.Range("K97").Value = 1897.0968
$ws.Range("L97").Value = 5640.1665
$ws.Range("M97").Value = -1401.0968
$ws.Range("N97").Value = -6632.1665
# row 122
$ws.Range("H122").Value = 1700
$ws.Range("I122").Value = 1800
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 5400
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -2950
$ws.Range("N122").Value = -9400

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 93
$ws.Range("H93").Value = 1122.7273
$ws.Range("I93").Value = 725
$ws.Range("J93").Value = 1211.1111
$ws.Range("K93").Value = 725
$ws.Range("L93").Value = 1211.1111
$ws.Range("M93").Value = 523
$ws.Range("N93").Value = -3707.1111
# row 122
$ws.Range("H122").Value = 113333.89
$ws.Range("I122").Value = 251876.25
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 755628.75
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -753178.75
$ws.Range("N122").Value = -12400

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 51
$ws.Range("H51").Value = 10000
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
# row 52
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
# row 122
$ws.Range("H122").Value = 9527809
$ws.Range("I122").Value = 14290714
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 42872142
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -42869692
$ws.Range("N122").Value = -10900

Write-Host "Edits applied successfully"
